$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 45: EEG Wind Energy Subsidy -------------------------------------
# Set values in the exact order the original workbook's shared-string table
# was built in, so newly created shared strings line up with the target.
$ws.Range("A45").Value = "eegWind"
$ws.Range("B45").Value = "EEG Wind Energy Subsidy"
$ws.Range("D45").Value = "Climate Policy"
$ws.Range("J45").Value = "https://ideas.repec.org/a/eee/pubeco/v169y2019icp172-202.html"
$ws.Range("K45").Value = "abrell2019"
$ws.Range("I45").Value = "Abrell et al. (2019)"

# --- Row 46: EEG Solar Energy Subsidy -------------------------------------
$ws.Range("A46").Value = "eegSolar"
$ws.Range("F45").Value = "In Germany electricity produced from renewable sources is subsidized under the the ""Erneuerbare Energien Gesetz"".  Abrell et al. (2019) estimate the implied carbon emission abatement cost of wind Energy."
$ws.Range("F46").Value = "In Germany electricity produced from renewable sources is subsidized under the the ""Erneuerbare Energien Gesetz"". Abrell et al. (2019) estimate the implied carbon emission abatement cost of solar Energy."
$ws.Range("B46").Value = "EEG Solar Energy Subsidy"

# --- Remaining (reused) values ---------------------------------------------
$ws.Range("C45").Value = 2012
$ws.Range("C46").Value = 2012
$ws.Range("D46").Value = "Climate Policy"
$ws.Range("I46").Value = "Abrell et al. (2019)"
$ws.Range("J46").Value = "https://ideas.repec.org/a/eee/pubeco/v169y2019icp172-202.html"
$ws.Range("K46").Value = "abrell2019"

# --- Styling: column F = wrap text -----------------------------------------
$ws.Range("F45").WrapText = $true
$ws.Range("F46").WrapText = $true

# --- Row heights -------------------------------------------------------------
$ws.Rows.Item(45).RowHeight = 75
$ws.Rows.Item(46).RowHeight = 75

# --- Hyperlinks (added before the "Link" style so the style is not
#     overwritten by Hyperlinks.Add's own default formatting) ---------------
$ws.Hyperlinks.Add($ws.Range("J45"), "https://ideas.repec.org/a/eee/pubeco/v169y2019icp172-202.html")
$ws.Hyperlinks.Add($ws.Range("J46"), "https://ideas.repec.org/a/eee/pubeco/v169y2019icp172-202.html")

# --- Styling: column J = hyperlink ("Link") style ---------------------------
$ws.Range("J45").Style = "Link"
$ws.Range("J46").Style = "Link"

# --- View state: keep tab selected, move viewport/selection near the bottom -
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 40
$win.ScrollColumn = 2
$ws.Range("B47").Select()
